# Refresh the crypto-ranking snapshot (price + 1h change columns,
# plus a newly-inserted "ApeXProtocol" row that pushes rows 45-51
# down by one and drops the old trailing "BitcoinSV" row off the
# bottom of the A1:E51 used range).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.914.53'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '2.301.33'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.511'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.57%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.504'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').Value = '2.659.47'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '2.301.49'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '42.852.01'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.60'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.05%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').Value = '2.013.92'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.63%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.07%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.526.88'
$ws.Range('E51').Value = '  -0.51%  '
